$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, shifting existing rows 129-143 down to 130-144.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly data record.
$ws.Cells.Item(129, 1).Value = 5
$ws.Cells.Item(129, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(129, 3).Value = "Maule"
$ws.Cells.Item(129, 4).Value = 44984
$ws.Cells.Item(129, 5).Value = 7
$ws.Cells.Item(129, 6).Value = 100112001
$ws.Cells.Item(129, 7).Value = "Berenjena"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 200
$ws.Cells.Item(129, 11).Value = 8000
$ws.Cells.Item(129, 12).Value = 8000
$ws.Cells.Item(129, 13).Value = 8000
$ws.Cells.Item(129, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(129, 15).Value = "Región del Maule"
$ws.Cells.Item(129, 16).Value = 160
$ws.Cells.Item(129, 17).Value = 50
$ws.Cells.Item(129, 18).Value = "Hortaliza"
